$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row additions (I1:AB1)
$headers = @("bfsp","bfrp","rest0","rest1","rest2","rest3","error","rnge","hold","wp","balk","pickoff","durability","dldays","bat","bunt","run","steal","ab","wild")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 9 + $i   # column I = 9
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# New data row additions (I2:AB2)
$ws.Cells.Item(2, 9).Value  = 0          # I2
$ws.Cells.Item(2, 10).Value = 6          # J2
$ws.Cells.Item(2, 11).Value = 11         # K2
$ws.Cells.Item(2, 12).Value = 14         # L2
$ws.Cells.Item(2, 13).Value = 17         # M2
$ws.Cells.Item(2, 14).Value = 19         # N2
$ws.Cells.Item(2, 15).Value = 4          # O2
$ws.Cells.Item(2, 16).Value = "d"        # P2
$ws.Cells.Item(2, 17).Value = "Vg_0"     # Q2
$ws.Cells.Item(2, 18).Value = "unlikely" # R2
$ws.Cells.Item(2, 19).Value = "rare"     # S2
$ws.Cells.Item(2, 20).Value = 0          # T2
$ws.Cells.Item(2, 21).Value = 7          # U2
$ws.Cells.Item(2, 22).Value = 0          # V2
$ws.Cells.Item(2, 23).Value = "rsp#1_pr" # W2
$ws.Cells.Item(2, 24).Value = "Fr_-1"    # X2
$ws.Cells.Item(2, 25).Value = 0          # Y2
$ws.Cells.Item(2, 26).Value = 0          # Z2
$ws.Cells.Item(2, 27).Value = 0          # AA2
$ws.Cells.Item(2, 28).Value = "null"     # AB2

# Update sheet view to match new active cell / scroll position
$ws.Application.ActiveWindow.ScrollColumn = 22  # topLeftCell = V1
$ws.Range("AB2").Select() | Out-Null
